$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must remain stored as TEXT even when it is
# numeric-looking (e.g. "0.242"), without leaving any NumberFormat/style
# residue behind. We build it as a text formula ("=""0.242""") and then
# convert it in-place to a literal value via Copy / PasteSpecial(values),
# which collapses the formula to a plain shared-string cell.
function Set-TextValue {
    param($addr, $text)
    $r = $ws.Range($addr)
    $escaped = $text -replace '"', '""'
    $r.Formula = '="' + $escaped + '"'
    $r.Copy() | Out-Null
    $r.PasteSpecial(-4163) | Out-Null
}

# --- Row 1 header labels: strip " Diff-in-Diff" suffix ---
$ws.Range("B1").Value2 = "C"
$ws.Range("C1").Value2 = "U"
$ws.Range("D1").Value2 = '$\pi$'
$ws.Range("E1").Value2 = "FFR"
$ws.Range("F1").Value2 = "A"

# --- Column A row labels: strip " Diff-in-Diff" suffix ---
$ws.Range("A2").Value2 = "C"
$ws.Range("A3").Value2 = "U"
$ws.Range("A4").Value2 = '$\pi$'
$ws.Range("A5").Value2 = "FFR"
$ws.Range("A6").Value2 = "A"

# --- Updated coefficient text values (column B) ---
Set-TextValue "B3" "0.094***"
Set-TextValue "B4" "-1.013**"
Set-TextValue "B5" "0.242"
Set-TextValue "B6" "4.17***"
Set-TextValue "B7" "-0.054"

# --- Updated coefficient text values (column C) ---
Set-TextValue "C2" "0.858***"
Set-TextValue "C4" "-3.322**"
Set-TextValue "C5" "-0.052"
Set-TextValue "C6" "-4.73**"
Set-TextValue "C7" "0.092"

# --- Updated coefficient text values (column D) ---
Set-TextValue "D2" "-0.047**"
Set-TextValue "D5" "-0.031"
Set-TextValue "D6" "-0.586***"
Set-TextValue "D7" "-0.016"

# --- Updated coefficient text values (column E) ---
Set-TextValue "E2" "0.037"
Set-TextValue "E4" "-0.101"
Set-TextValue "E6" "-1.02***"
Set-TextValue "E7" "-0.124"

# --- Updated coefficient text values (column F) ---
Set-TextValue "F4" "-0.264***"
Set-TextValue "F5" "-0.138***"

# --- r2_adj row: plain numeric value change ---
$ws.Range("D8").Value2 = 0.57
